$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header columns (row 1)
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# Title-case "de/del/el/la/las/los/y" connector words in state/municipality names
$renames = @{
  "B6" = "Pabellón De Arteaga"
  "B7" = "Rincón De Romos"
  "B22" = "Amatenango De La Frontera"
  "B32" = "Comitán De Domínguez"
  "B45" = "Mazapa De Madero"
  "B52" = "Salto De Agua"
  "B53" = "San Cristóbal De Las Casas"
  "B81" = "Guadalupe Y Calvo"
  "B83" = "Hidalgo Del Parral"
  "B89" = "San Francisco Del Oro"
  "B105" = "San Juan De Sabinas"
  "A113" = "Ciudad De México"
  "B117" = "Cuajimalpa De Morelos"
  "B131" = "Coneto De Comonfort"
  "B146" = "Pánuco De Coronado"
  "B149" = "San Juan Del Río"
  "A154" = "Estado De México"
  "B154" = "Acambay De Ruíz Castañeda"
  "B156" = "Almoloya De Alquisiras"
  "B161" = "Atizapán De Zaragoza"
  "B165" = "Chapa De Mota"
  "B168" = "Coacalco De Berriozábal"
  "B172" = "Ecatepec De Morelos"
  "B176" = "Ixtapan De La Sal"
  "B185" = "Naucalpan De Juárez"
  "B193" = "San Felipe Del Progreso"
  "B201" = "Tenango Del Valle"
  "B207" = "Tlalnepantla De Baz"
  "B211" = "Valle De Bravo"
  "B212" = "Villa De Allende"
  "B221" = "San Miguel De Allende"
  "B222" = "Apaseo El Alto"
  "B223" = "Apaseo El Grande"
  "B230" = "Dolores Hidalgo Cuna De La Independencia Nacional"
  "B234" = "Jaral Del Progreso"
  "B241" = "Purísima Del Rincón"
  "B245" = "San Diego De La Unión"
  "B247" = "San Francisco Del Rincón"
  "B249" = "San Luis De La Paz"
  "B251" = "Silao De La Victoria"
  "B256" = "Valle De Santiago"
  "B262" = "Acapulco De Juárez"
  "B264" = "Ajuchitlán Del Progreso"
  "B265" = "Alcozauca De Guerrero"
  "B269" = "Atenango Del Río"
  "B270" = "Atlamajalcingo Del Monte"
  "B272" = "Atoyac De Álvarez"
  "B273" = "Ayutla De Los Libres"
  "B276" = "Chilapa De Álvarez"
  "B277" = "Chilpancingo De Los Bravo"
  "B281" = "Coyuca De Benítez"
  "B282" = "Coyuca De Catalán"
  "B284" = "Cutzamala De Pinzón"
  "B290" = "Huitzuco De Los Figueroa"
  "B291" = "Iguala De La Independencia"
  "B292" = "Ixcateopan De Cuauhtémoc"
  "B295" = "La Unión De Isidoro Montes De Oca"
  "B298" = "Mártir De Cuilapan"
  "B311" = "Taxco De Alarcón"
  "B313" = "Técpan De Galeana"
  "B315" = "Tepecoacuilco De Trujano"
  "B316" = "Tixtla De Guerrero"
  "B319" = "Tlalixtaquilla De Maldonado"
  "B320" = "Tlapa De Comonfort"
  "B330" = "Agua Blanca De Iturbide"
  "B333" = "Atotonilco El Grande"
  "B338" = "Cuautepec De Hinojosa"
  "B342" = "Huejutla De Reyes"
  "B345" = "Jacala De Ledezma"
  "B350" = "Mineral Del Chico"
  "B351" = "Mineral Del Monte"
  "B352" = "Mixquiahuala De Juárez"
  "B354" = "Pachuca De Soto"
  "B360" = "Santiago De Anaya"
  "B363" = "Tepehuacán De Guerrero"
  "B364" = "Tepeji Del Río De Ocampo"
  "B369" = "Tula De Allende"
  "B370" = "Tulancingo De Bravo"
  "B372" = "Zacualtipán De Ángeles"
  "B376" = "Ahualulco De Mercado"
  "B379" = "Atemajac De Brizuela"
  "B382" = "Atotonilco El Alto"
  "B391" = "Encarnación De Díaz"
  "B396" = "Huejuquilla El Alto"
  "B397" = "Ixtlahuacán Del Río"
  "B400" = "Jilotlán De Los Dolores"
  "B403" = "Lagos De Moreno"
  "B409" = "Ojuelos De Jalisco"
  "B413" = "San Cristóbal De La Barranca"
  "B415" = "San Juan De Los Lagos"
  "B416" = "San Juanito De Escobedo"
  "B420" = "San Miguel El Alto"
  "B421" = "San Sebastián Del Oeste"
  "B424" = "Talpa De Allende"
  "B425" = "Tamazula De Gordiano"
  "B427" = "Tepatitlán De Morelos"
  "B428" = "Tizapán El Alto"
  "B429" = "Tlajomulco De Zúñiga"
  "B435" = "Unión De San Antonio"
  "B436" = "Valle De Juárez"
  "B439" = "Yahualica De González Gallo"
  "B440" = "Zacoalco De Torres"
  "B443" = "Zapotlán Del Rey"
  "B444" = "Zapotlán El Grande"
  "B519" = "Tiquicheo De Nicolás Romero"
  "B550" = "Puente De Ixtla"
  "B558" = "Zacualpan De Amilpas"
  "B560" = "Bahía De Banderas"
  "B562" = "Ixtlán Del Río"
  "B587" = "San Nicolás De Los Garza"
  "B592" = "Acatlán De Pérez Figueroa"
  "B595" = "Chalcatongo De Hidalgo"
  "B596" = "Ciénega De Zimatlán"
  "B598" = "Fresnillo De Trujano"
  "B599" = "Heroica Ciudad De Ejutla De Crespo"
  "B600" = "Heroica Ciudad De Huajuapan De León"
  "B602" = "Ixtlán De Juárez"
  "B603" = "Heroica Ciudad De Juchitán De Zaragoza"
  "B607" = "Mariscala De Juárez"
  "B609" = "Miahuatlán De Porfirio Díaz"
  "B611" = "Nejapa De Madero"
  "B612" = "Oaxaca De Juárez"
  "B613" = "Ocotlán De Morelos"
  "B615" = "Putla Villa De Guerrero"
  "B626" = "San Felipe Jalapa De Díaz"
  "B637" = "San Juan Del Río"
  "B651" = "San Miguel Del Puerto"
  "B674" = "Santa María Jalapa Del Marqués"
  "B684" = "Santiago Del Río"
  "B699" = "Tamazulápam Del Espíritu Santo"
  "B701" = "Teotitlán De Flores Magón"
  "B704" = "Tlacolula De Matamoros"
  "B705" = "Totontepec Villa De Morelos"
  "B707" = "Villa De Tututepec De Melchor Ocampo"
  "B708" = "Villa De Zaachila"
  "B709" = "Villa Sola De Vega"
  "B710" = "Zimatlán De Álvarez"
  "B723" = "Ayotoxco De Guerrero"
  "B728" = "Chila De La Sal"
  "B732" = "Cuayuca De Andrade"
  "B733" = "Cuetzalan Del Progreso"
  "B744" = "Ixcamilpa De Guerrero"
  "B746" = "Izúcar De Matamoros"
  "B751" = "Los Reyes De Juárez"
  "B766" = "San Salvador El Verde"
  "B774" = "Tepatlaxco De Hidalgo"
  "B777" = "Tepexi De Rodríguez"
  "B779" = "Tetela De Ocampo"
  "B792" = "Xayacatlán De Bravo"
  "B801" = "Amealco De Bonfil"
  "B803" = "Cadereyta De Montes"
  "B807" = "Jalpan De Serra"
  "B808" = "Landa De Matamoros"
  "B810" = "Pinal De Amoles"
  "B812" = "San Juan Del Río"
  "B826" = "Ciudad Del Maíz"
  "B834" = "Mexquitic De Carmona"
  "B839" = "San Ciro De Acosta"
  "B843" = "Santa María Del Río"
  "B844" = "Soledad De Graciano Sánchez"
  "B850" = "Tanquián De Escobedo"
  "B852" = "Villa De Arriaga"
  "B853" = "Villa De Guadalupe"
  "B854" = "Villa De Ramos"
  "B855" = "Villa De Reyes"
  "B889" = "Jalpa De Méndez"
  "B925" = "Muñoz De Domingo Arenas"
  "B926" = "Nanacamilpa De Mariano Arista"
  "B940" = "Alto Lucero De Gutiérrez Barrios"
  "B942" = "Amatlán De Los Reyes"
  "B947" = "Camarón De Tejeda"
  "B949" = "Castillo De Teayo"
  "B957" = "Cosamaloapan De Carpio"
  "B967" = "Hueyapan De Ocampo"
  "B968" = "Ignacio De La Llave"
  "B969" = "Ixhuatlán De Madero"
  "B974" = "Juchique De Ferrer"
  "B977" = "Martínez De La Torre"
  "B988" = "Paso De Ovejas"
  "B989" = "Paso Del Macho"
  "B991" = "Poza Rica De Hidalgo"
  "B996" = "Sayula De Alemán"
  "B998" = "Soledad De Doblado"
  "B1002" = "Tatahuicapan De Juárez"
  "B1018" = "Vega De Alatorre"
  "B1035" = "Concepción Del Oro"
  "B1036" = "El Plateado De Joaquín Amaro"
  "B1050" = "Moyahua De Estrada"
  "B1051" = "Noria De Ángeles"
  "B1059" = "Teúl De González Ortega"
  "B1060" = "Tlaltenango De Sánchez Román"
  "B1063" = "Villa De Cos"
}

foreach ($key in $renames.Keys) {
  $ws.Range($key).Value = $renames[$key]
}

# Remove trailing metadata/footnote rows (1072-1076); row 1070 (grand Total) stays, row 1071 was already empty
$ws.Range("A1072:A1076").EntireRow.Delete() | Out-Null

Write-Output ("done; dim=" + $ws.UsedRange.Address())